$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The extrapolation calculation for column C (HDP kraje na obyvatele) was
# changed - all values are scaled down by a constant factor
# (~0.9708179410257477) to account for the newly added RUD-per-capita
# column's effect on the underlying computation.
$newValues = @{
    2  = 636733.2941520028
    3  = 622606.6389124633
    4  = 445948.2001064043
    5  = 632173.8398814903
    6  = 645430.2588331797
    7  = 576964.7709690979
    8  = 600611.1212153231
    9  = 535593.3298955993
    10 = 541483.8851404014
    11 = 656266.418216069
    12 = 613106.0944495659
    13 = 578096.3282618701
    14 = 712858.6685511762
}

foreach ($row in $newValues.Keys) {
    $ws.Range("C$row").Value = $newValues[$row]
}
